$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.990.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "'1.555.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").Value = "'207.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "'0.487"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.248"
$ws.Range("D8").Style = "Normal"

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'21.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").Value = "'0.0589"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("D12").Value = "'1.778.56"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'1.557.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("E14").Value = "  +1.33%  "

$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'26.993.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'61.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").Value = "'215.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").Value = "'4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("E23").Value = "  +2.92%  "

$ws.Range("E24").Value = "  -0.94%  "

$ws.Range("D25").Value = "'152.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("E29").Value = "  +1.45%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").Value = "'1.403.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.74%  "

$ws.Range("E34").Value = "  +3.04%  "

$ws.Range("E35").Value = "  +3.40%  "

$ws.Range("D36").Value = "'0.951"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.66%  "

$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("E39").Value = "  +0.58%  "

$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("D42").Value = "'0.987"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.71%  "

$ws.Range("E43").Value = "  +3.68%  "

$ws.Range("D44").Value = "'5.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.18%  "

$ws.Range("D45").Value = "'63.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "

$ws.Range("D46").Value = "'1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").Value = "'1.691.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "'86.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").Value = "'0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("E51").Value = "  +0.66%  "
